$d = $word.ActiveDocument

$newText = "Waktu Kampanye Orion: 16-25 Januari, 14-23 Februari, 14-24 Maret"
$oldText = "Waktu Kampanye 2018 untuk Perseus: 30 Oktober-8 November dan 29 November-8 Desember"

function Replace-WithPlainRun($searchText) {
    $rng = $word.ActiveDocument.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $startPos = $rng.Start
        # Remove all matched (possibly multi-run) text, collapsing the paragraph's
        # run content so the subsequent insertion starts with a clean slate (no
        # inherited rPr from any neighboring run).
        $rng.Text = ""
        $insPoint = $word.ActiveDocument.Range($startPos, $startPos)
        $insPoint.InsertAfter($newText)
    }
    return $found
}

# Occurrence 1: paragraph begins with a lone red space run, then "Waktu Kampanye "
# then "2018 untuk Perseus...". All three runs collapse into one plain run.
$found1 = Replace-WithPlainRun(" " + $oldText)
Write-Output "occurrence 1 replaced: $found1"

# Occurrence 2: paragraph begins with "www.globeatnight.org" + line break, then
# "Waktu Kampanye " then "2018 untuk Perseus...". All runs collapse into one.
$found2 = Replace-WithPlainRun("www.globeatnight.org" + [char]11 + $oldText)
Write-Output "occurrence 2 replaced: $found2"

# Occurrences 3 and 4: paragraph begins directly with "Waktu Kampanye " then
# "2018 untuk Perseus...". Both runs collapse into one. Two such paragraphs remain.
$found3 = Replace-WithPlainRun($oldText)
Write-Output "occurrence 3 replaced: $found3"

$found4 = Replace-WithPlainRun($oldText)
Write-Output "occurrence 4 replaced: $found4"
